$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.903.25'
$ws.Range('E2').Value = '  -0.20%  '

$ws.Range('D3').Value = '1.635.42'
$ws.Range('E3').Value = '  -0.32%  '

$ws.Range('E4').Value = '  -0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.25%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5058'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.22%  '

$ws.Range('E7').Value = '  -0.19%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2577'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.73%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06359'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.13%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.62'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.63%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07742'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.37%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.269'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.28%  '

$ws.Range('D13').Value = '1.626.66'
$ws.Range('E13').Value = '  -0.89%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5469'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.41%  '

$ws.Range('D15').Value = '0.0₅7732'
$ws.Range('E15').Value = '  -1.09%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.14'
$ws.Range('D16').Style = 'Normal'

$ws.Range('D17').Value = '25.894.20'
$ws.Range('E17').Value = '  -0.32%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.440'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.14%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.83'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.30%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.895'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.53%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.074'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.51%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.002'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.27%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.901'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.18%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '142.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.67%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1245'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.16%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.824'
$ws.Range('D27').Style = 'Normal'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.56'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.97%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.241'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.35%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.04864'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.67%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.238'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.54%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.192'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.17%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.544'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.36%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.367'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.27%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9065'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.59%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.574'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.38%  '

$ws.Range('D37').Value = '1.126.59'
$ws.Range('E37').Value = '  -0.21%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5494'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.06%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01556'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.02%  '

$ws.Range('E40').Value = '  -0.20%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.591'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8024'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.70%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '98.23'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.55%  '

$ws.Range('E44').Value = '  -5.30%  '

$ws.Range('D45').Value = '1.769.78'
$ws.Range('E45').Value = '  -0.37%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4461'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.74%  '

$ws.Range('E47').Value = '  -0.21%  '

$ws.Range('E48').Value = '  +0.17%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05157'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.65%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.524'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.20%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.003'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.05%  '

